# Add the "Special Character Removal" sheet to the end of the workbook
# (after "Rename Column") and fill it with the same Action/Time/Content
# table pattern used by the other process sheets.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Special Character Removal"

# -- Header row --
$ws.Range("A1").Value = "Action"
$ws.Range("B1").Value = "Time"
$ws.Range("C1").Value = "Content"

# -- Row 2 --
$ws.Range("A2").Value = "Upload CSV"
$ws.Range("B2").Value = "5 min"
$ws.Range("C2").Value = "df = pd.read_csv('file.csv')"

# -- Row 3 --
$ws.Range("A3").Value = "Identify Characters"
$ws.Range("B3").Value = "2 min"
$ws.Range("C3").Value = "Visual inspection with df.head()"

# -- Row 4 --
$ws.Range("A4").Value = "Remove Characters"
$ws.Range("B4").Value = "2 min"
$ws.Range("C4").Value = "df.replace({r'[^\x00-\x7F]+':''}, regex=True, inplace=True)"

# -- Row 5 --
$ws.Range("A5").Value = "Verify Changes"
$ws.Range("B5").Value = "1 min"
$ws.Range("C5").Value = "df.head() to check cleaned data"

# -- Row 6 (summary) --
$ws.Range("A6").Value = "Overall"
$ws.Range("B6").Value = "10 min"

# -- Formatting, matching the other sheets (bold 13pt header/summary rows,
#    regular 13pt body rows, ~17pt row height) --
$ws.Range("A1:C5").Font.Size = 13
$ws.Range("A6:B6").Font.Size = 13
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A6:B6").Font.Bold = $true

$ws.Range("A1:C5").RowHeight = 17
$ws.Range("A6:B6").RowHeight = 17

# Leave the new sheet as the active one, with the used range selected.
[void]$ws.Range("A1:C6").Select()
